$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the data row currently at row 2 (task #48) down to become row 49,
# shifting the rows currently at 3..49 up by one (to 2..48). Row 50
# onward stays where it is.
$movedRow = $ws.Range("A2:F2").Value()
$block = $ws.Range("A3:F49").Value()

$ws.Range("A2:F48").Value2 = $block
$ws.Range("A49:F49").Value2 = $movedRow

# Update the view state: select row 2 (matches the new top of the list)
$ws.Rows.Item(2).Select() | Out-Null
